$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add 5 new rows (7-11) mirroring rows 2-6 but referencing Table[2] xpaths
$ws.Range("B7").Value = "/Envelope/Body/GetInfoByAreaCodeResponse/GetInfoByAreaCodeResult/NewDataSet/Table[2]/CITY"
$ws.Range("C7").Value = "[A-Z a-z].*"

$ws.Range("B8").Value = "/Envelope/Body/GetInfoByAreaCodeResponse/GetInfoByAreaCodeResult/NewDataSet/Table[2]/STATE"
$ws.Range("C8").Value = "[A-Z]{2}"

$ws.Range("B9").Value = "/Envelope/Body/GetInfoByAreaCodeResponse/GetInfoByAreaCodeResult/NewDataSet/Table[2]/ZIP"
$ws.Range("C9").Value = "[0-9]{5}"

$ws.Range("B10").Value = "/Envelope/Body/GetInfoByAreaCodeResponse/GetInfoByAreaCodeResult/NewDataSet/Table[2]/AREA_CODE"
$ws.Range("C10").Value = "[0-9]{3}"

$ws.Range("B11").Value = "/Envelope/Body/GetInfoByAreaCodeResponse/GetInfoByAreaCodeResult/NewDataSet/Table[2]/TIME_ZONE"
$ws.Range("C11").Value = "[A-Z]{1}"

# Update selection to match the target state (active cell B11)
$ws.Range("B11").Select()
